$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'37.558.08"
$ws.Range("E2").Value = "  +5.22%  "

# Row 3
$ws.Range("D3").Value = "'2.059.70"
$ws.Range("E3").Value = "  +3.84%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'254.71"
$ws.Range("E5").Value = "  +3.65%  "

# Row 6
$ws.Range("E6").Value = "  +3.13%  "

# Row 7
$ws.Range("D7").Value = "'66.70"
$ws.Range("E7").Value = "  +12.14%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("E9").Value = "  +9.17%  "

# Row 10
$ws.Range("D10").Value = "'60.32"
$ws.Range("E10").Value = "  +2.32%  "

# Row 11
$ws.Range("D11").Value = "'0.0775"
$ws.Range("E11").Value = "  +4.73%  "

# Row 12
$ws.Range("E12").Value = "  +0.54%  "

# Row 13
$ws.Range("D13").Value = "'0.943"
$ws.Range("E13").Value = "  -0.77%  "

# Row 14
$ws.Range("D14").Value = "'24.04"
$ws.Range("E14").Value = "  +27.46%  "

# Row 15
$ws.Range("D15").Value = "'15.18"
$ws.Range("E15").Value = "  +3.57%  "

# Row 16
$ws.Range("D16").Value = "'2.360.50"
$ws.Range("E16").Value = "  +3.86%  "

# Row 17
$ws.Range("E17").Value = "  +8.26%  "

# Row 18
$ws.Range("D18").Value = "'2.057.31"
$ws.Range("E18").Value = "  +3.82%  "

# Row 19
$ws.Range("D19").Value = "'37.472.62"
$ws.Range("E19").Value = "  +5.29%  "

# Row 20
$ws.Range("E20").Value = "  +3.00%  "

# Row 21
$ws.Range("D21").Value = "'0.0₃0883"
$ws.Range("E21").Value = "  +3.94%  "

# Row 22
$ws.Range("E22").Value = "  +5.77%  "

# Row 23
$ws.Range("D23").Value = "'241.44"
$ws.Range("E23").Value = "  +3.47%  "

# Row 24
$ws.Range("D24").Value = "'2.67"
$ws.Range("E24").Value = "  +4.65%  "

# Row 25
$ws.Range("E25").Value = "  -0.10%  "

# Row 26
$ws.Range("E26").Value = "  +7.54%  "

# Row 27
$ws.Range("D27").Value = "'10.11"
$ws.Range("E27").Value = "  +9.77%  "

# Row 28
$ws.Range("D28").Value = "'162.61"
$ws.Range("E28").Value = "  -1.63%  "

# Row 29
$ws.Range("D29").Value = "'0.137"
$ws.Range("E29").Value = "  +42.92%  "

# Row 30
$ws.Range("D30").Value = "'20.24"
$ws.Range("E30").Value = "  +4.87%  "

# Row 31
$ws.Range("E31").Value = "  +3.35%  "

# Row 32
$ws.Range("E32").Value = "  +6.98%  "

# Row 33
$ws.Range("E33").Value = "  +8.17%  "

# Row 34
$ws.Range("D34").Value = "'0.0639"
$ws.Range("E34").Value = "  +6.83%  "

# Row 35
$ws.Range("D35").Value = "'4.72"
$ws.Range("E35").Value = "  +8.21%  "

# Row 36
$ws.Range("D36").Value = "'2.46"
$ws.Range("E36").Value = "  +0.33%  "

# Row 37
$ws.Range("D37").Value = "'6.34"
$ws.Range("E37").Value = "  +15.28%  "

# Row 38
$ws.Range("E38").Value = "  -0.03%  "

# Row 39
$ws.Range("E39").Value = "  +2.32%  "

# Row 40
$ws.Range("E40").Value = "  +38.48%  "

# Row 41
$ws.Range("D41").Value = "'0.105"
$ws.Range("E41").Value = "  +13.93%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.28"
$ws.Range("E42").Value = "  +4.69%  "

# Row 43
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'18.28"
$ws.Range("E43").Value = "  +10.85%  "

# Row 44
$ws.Range("E44").Value = "  +6.95%  "

# Row 45
$ws.Range("D45").Value = "'1.17"
$ws.Range("E45").Value = "  +6.88%  "

# Row 46
$ws.Range("E46").Value = "  +3.29%  "

# Row 47
$ws.Range("D47").Value = "'97.46"
$ws.Range("E47").Value = "  +4.03%  "

# Row 48
$ws.Range("D48").Value = "'8.03"
$ws.Range("E48").Value = "  +2.80%  "

# Row 49
$ws.Range("D49").Value = "'1.420.19"
$ws.Range("E49").Value = "  +4.02%  "

# Row 50
$ws.Range("E50").Value = "  +2.21%  "

# Row 51
$ws.Range("D51").Value = "'48.60"
$ws.Range("E51").Value = "  +3.21%  "
